$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.647599999999994
$ws.Range("B8").Value = 4.613600000000002
$ws.Range("B10").Value = 8.547700000000006
$ws.Range("B12").Value = 5.5959
$ws.Range("C13").Value = -12.4118
$ws.Range("B18").Value = 4.917800000000003
$ws.Range("D20").Value = -8.305300000000004
$ws.Range("B25").Value = 5.964199999999995
